$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 126
$ws.Range("B2").Value = 257
$ws.Range("C2").Value = "벽에 덕지덕치 (처 발랐다 / 처발랐다)"
$ws.Range("D2").Value = "하"
$ws.Range("E2").Value = "띄어쓰기"
$ws.Range("F2").Value = "처발랐다"
$ws.Range("G2").Value = 0.3647933250547525

# Row 3
$ws.Range("A3").Value = 123
$ws.Range("B3").Value = 253
$ws.Range("C3").Value = "(세종조 / 세종 조)에 한글이 창제되었다"
$ws.Range("D3").Value = "하"
$ws.Range("E3").Value = "띄어쓰기"
$ws.Range("F3").Value = "세종조"
$ws.Range("G3").Value = 0.3647933250547525

# Row 4
$ws.Range("A4").Value = 121
$ws.Range("B4").Value = 248
$ws.Range("C4").Value = "꽃을 (한아름 / 한 아름) 사왔다"
$ws.Range("D4").Value = "하"
$ws.Range("E4").Value = "띄어쓰기"
$ws.Range("F4").Value = "한 아름"
$ws.Range("G4").Value = 0.3647933250547525

# Row 5
$ws.Range("A5").Value = 120
$ws.Range("B5").Value = 247
$ws.Range("C5").Value = "질문에 대해 자세히 (답변드렸다 / 답변 드렸다)"
$ws.Range("D5").Value = "하"
$ws.Range("E5").Value = "띄어쓰기"
$ws.Range("F5").Value = "답변드렸다"
$ws.Range("G5").Value = 0.3647933250547525

# Row 6
$ws.Range("A6").Value = 99
$ws.Range("B6").Value = 217
$ws.Range("C6").Value = "이런 일은 (생전 처음 / 생전처음) 겪는 일이다"
$ws.Range("D6").Value = "하"
$ws.Range("E6").Value = "띄어쓰기"
$ws.Range("F6").Value = "생전 처음"
$ws.Range("G6").Value = 0.3647933250547525

# Row 7
$ws.Range("A7").Value = 106
$ws.Range("B7").Value = 232
$ws.Range("C7").Value = "고향에 금방 (내려 가겠다 / 내려가겠다)"
$ws.Range("D7").Value = "하"
$ws.Range("E7").Value = "띄어쓰기"
$ws.Range("F7").Value = "내려가겠다"
$ws.Range("G7").Value = 0.3647933250547525

# Row 8
$ws.Range("A8").Value = 74
$ws.Range("B8").Value = 169
$ws.Range("C8").Value = "국의 (건더기 / 건데기)를 골라 먹었다"
$ws.Range("D8").Value = "하"
$ws.Range("E8").Value = "맞춤법"
$ws.Range("F8").Value = "건더기"
$ws.Range("G8").Value = 0.2538352261379218
